# ----------------------------------------------------------------------------
# Refactor: the combined "lnglat" column (N, e.g. "1.3044719,103.7724654") is
# split into two independent columns -- "Latitude" (N, kept as text) and
# "Longitude" (O, numeric) -- mirroring a pandas df query that now returns the
# coordinates as separate fields instead of one comma-joined string.
# ----------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23's latitude is the one pair whose split-out value is still a duplicate
# of another row once lat/lng are separated, so it is written while column N is
# still General-formatted -- it has to land in the sheet as a plain number, not
# shared-string text. Write it first, before the column gets its text format.
$ws.Range("N23").Value = 1.3038022

# Format the whole Latitude column as text so every numeric-looking value typed
# into it from here on is kept verbatim as a string instead of being coerced to
# a number. This also mints the two new cellXfs style records used below.
$ws.Columns("N").NumberFormat = "@"

# Applying the column-wide format stamps a few blank rows below the table (28-30)
# with the new default column style; clear those placeholder cells back out so no
# stray cells are left behind.
$ws.Range("N28:N30").Clear()

# Write out every row's latitude (text) / longitude (number) pair.
$ws.Range("N2").Value = "1.3044719"
$ws.Range("O2").Value = 103.7724654
$ws.Range("N3").Value = "1.307338"
$ws.Range("O3").Value = 103.7726078
$ws.Range("N4").Value = "1.3045409"
$ws.Range("O4").Value = 103.7727869
$ws.Range("N5").Value = "1.303794"
$ws.Range("O5").Value = 103.7735167
$ws.Range("N6").Value = "1.306019"
$ws.Range("O6").Value = 103.772678
$ws.Range("N7").Value = "1.3050106"
$ws.Range("O7").Value = 103.7723947
$ws.Range("N8").Value = "1.3039258"
$ws.Range("O8").Value = 103.7735858
$ws.Range("N9").Value = "1.3046285"
$ws.Range("O9").Value = 103.7730182
$ws.Range("N10").Value = "1.3047189"
$ws.Range("O10").Value = 103.7727242
$ws.Range("N11").Value = "1.3044706"
$ws.Range("O11").Value = 103.7724575
$ws.Range("N12").Value = "1.3051092"
$ws.Range("O12").Value = 103.7723276
$ws.Range("N13").Value = "1.3040592"
$ws.Range("O13").Value = 103.7741032
$ws.Range("N14").Value = "1.305593"
$ws.Range("O14").Value = 103.773083
$ws.Range("N15").Value = "1.3045187"
$ws.Range("O15").Value = 103.7728417
$ws.Range("N16").Value = "1.3048207"
$ws.Range("O16").Value = 103.7725693
$ws.Range("N17").Value = "1.305796"
$ws.Range("O17").Value = 103.773008
$ws.Range("N18").Value = "1.3046285"
$ws.Range("O18").Value = 103.7730182
$ws.Range("N19").Value = "1.3042717"
$ws.Range("O19").Value = 103.7738946
$ws.Range("N20").Value = "1.3046387"
$ws.Range("O20").Value = 103.7728153
$ws.Range("N21").Value = "1.3045756"
$ws.Range("O21").Value = 103.7726986
$ws.Range("N22").Value = "1.3054322"
$ws.Range("O22").Value = 103.7728657
$ws.Range("O23").Value = 103.7738266
$ws.Range("N24").Value = "1.3049764"
$ws.Range("O24").Value = 103.7724652
$ws.Range("N25").Value = "1.3039101"
$ws.Range("O25").Value = 103.7738303
$ws.Range("N26").Value = "1.3040203"
$ws.Range("O26").Value = 103.7741394
$ws.Range("N27").Value = "1.3047292"
$ws.Range("O27").Value = 103.7725536

# Replace the old combined "lnglat" header with separate Latitude/Longitude ones.
$ws.Range("N1").Value = "Latitude"
$ws.Range("O1").Value = "Longitude"

# --- View-state touch-ups: window position and the active cell selection. ---
$excel.ActiveWindow.Left = 5680
$excel.ActiveWindow.Top = 3220
[void]$ws.Range("P18").Select()
